$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.920.56'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '2.586.40'
$ws.Range('E3').Value = '  +2.63%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.04'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.02%  '
$ws.Range('E7').Value = '  +1.00%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.551'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.63'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E11').Value = '  +1.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.63'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.06%  '
$ws.Range('E13').Value = '  +7.38%  '
$ws.Range('D14').Value = '2.562.82'
$ws.Range('E14').Value = '  +1.60%  '
$ws.Range('E15').Value = '  +3.07%  '
$ws.Range('E16').Value = '  +2.83%  '
$ws.Range('D17').Value = '43.021.28'
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.94'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.46%  '
$ws.Range('E19').Value = '  +3.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.67'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.03'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '255.05'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.32%  '
$ws.Range('E23').Value = '  +2.44%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '28.75'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.88%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '39.41'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.76%  '
$ws.Range('E29').Value = '  -0.62%  '
$ws.Range('E30').Value = '  +1.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '155.73'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.00%  '
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.76'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0814'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.08%  '
$ws.Range('E35').Value = '  -2.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.37'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +12.69%  '
$ws.Range('E37').Value = '  +0.74%  '
$ws.Range('E38').Value = '  +1.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '23.57'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.11%  '
$ws.Range('E40').Value = '  +2.46%  '
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0312'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.08'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +29.43%  '
$ws.Range('D44').Value = '2.069.37'
$ws.Range('E44').Value = '  +2.70%  '
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.25'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.54%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.36'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '76.90'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +12.85%  '
$ws.Range('D49').Value = '2.833.95'
$ws.Range('E49').Value = '  +2.65%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '106.28'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.57%  '
$ws.Range('E51').Value = '  +2.81%  '
